# Captain's Quest Treasure Island — add Meta description paragraph under the
# title, drop the duplicated "Play ..." heading-style paragraph near the end,
# and repurpose the trailing italic paragraph into the feature-image prompt.

$d = $word.ActiveDocument

$wNs = "http://schemas.openxmlformats.org/wordprocessingml/2006/main"
$rsq = [char]0x2019   # RIGHT SINGLE QUOTATION MARK  U+2019  ( ' )

# ---------------------------------------------------------------------------
# 1) Insert the new "Meta description" paragraph right after the Heading1
#    title paragraph ("Play Captain's Quest Treasure Island Free Slot Game").
# ---------------------------------------------------------------------------
$title = $d.Paragraphs(1)
$null = $title.Range.InsertParagraphAfter()
$metaPara = $d.Paragraphs(2)

$metaText = ": Read our review of Captain" + $rsq + "s Quest Treasure Island, the pirate-themed slot game. Play now for free and experience the high volatility, Free Spins bonus rounds, and immersive graphics."
$metaXml = "<w:p xmlns:w='$wNs'><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Meta description</w:t></w:r><w:r><w:t>$metaText</w:t></w:r></w:p>"
$null = $metaPara.Range.InsertXML($metaXml)

# ---------------------------------------------------------------------------
# 2) Remove the trailing bold "Play Captain's Quest Treasure Island Free Slot
#    Game" paragraph, and turn the final italic paragraph into the new
#    feature-image art-direction prompt (keeping its italic formatting).
# ---------------------------------------------------------------------------
$total = $d.Paragraphs.Count
$boldPromoPara = $d.Paragraphs($total - 1)
$null = $boldPromoPara.Range.Delete()

$total = $d.Paragraphs.Count
$imagePromptPara = $d.Paragraphs($total)

$imageText = "For the feature image, create a cartoon-style design featuring a Maya warrior with glasses who is looking happy and satisfied. The design should include elements of the game, such as a ship sailing the Caribbean Sea, a deserted island where the treasure is hidden, and symbols of the game like the poker card suits, the helm, and the treasure. The background of the image should be blue with a pirate-themed border, and the game's name `"Captain's Quest Treasure Island`" should be prominently displayed. Make sure the image is bright and eye-catching, with lots of detail to entice players to try out the game."
$imageXml = "<w:p xmlns:w='$wNs'><w:r/><w:r><w:rPr><w:i/></w:rPr><w:t>$imageText</w:t></w:r></w:p>"
$null = $imagePromptPara.Range.InsertXML($imageXml)

Write-Output "done"
